# Convert the "Year" column from textual "NNNN BC" labels to negative
# numeric years (so the timeline can be sorted/plotted numerically), and
# move the active selection from A11 to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> numeric year (BC dates become negative numbers)
$years = @{
    2  = -8000
    3  = -7900
    4  = -7500
    5  = -7500
    6  = -7000
    7  = -6200
    8  = -5500
    9  = -5000
    10 = -4000
}

foreach ($row in $years.Keys) {
    $ws.Cells.Item($row, 1).Value = $years[$row]
}

[void]$ws.Range("C10").Select()
